$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column - copy formatting from E1 (adjacent header)
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:41:32.109076",
    "2021-10-05 13:41:32.109087",
    "2021-10-05 13:41:32.109090",
    "2021-10-05 13:41:32.109093",
    "2021-10-05 13:41:32.109097",
    "2021-10-05 13:41:32.109100",
    "2021-10-05 13:41:32.109102",
    "2021-10-05 13:41:32.109105",
    "2021-10-05 13:41:32.109108",
    "2021-10-05 13:41:32.109111",
    "2021-10-05 13:41:32.109114",
    "2021-10-05 13:41:32.109117",
    "2021-10-05 13:41:32.109120",
    "2021-10-05 13:41:32.109123",
    "2021-10-05 13:41:32.109126",
    "2021-10-05 13:41:32.109129",
    "2021-10-05 13:41:32.109132",
    "2021-10-05 13:41:32.109135",
    "2021-10-05 13:41:32.109138",
    "2021-10-05 13:41:32.109140",
    "2021-10-05 13:41:32.109143",
    "2021-10-05 13:41:32.109146",
    "2021-10-05 13:41:32.109149",
    "2021-10-05 13:41:32.109152",
    "2021-10-05 13:41:32.109155",
    "2021-10-05 13:41:32.109158",
    "2021-10-05 13:41:32.109160",
    "2021-10-05 13:41:32.109163",
    "2021-10-05 13:41:32.109166",
    "2021-10-05 13:41:32.109169",
    "2021-10-05 13:41:32.109171",
    "2021-10-05 13:41:32.109174",
    "2021-10-05 13:41:32.109177",
    "2021-10-05 13:41:32.109180",
    "2021-10-05 13:41:32.109183",
    "2021-10-05 13:41:32.109186",
    "2021-10-05 13:41:32.109189",
    "2021-10-05 13:41:32.109191",
    "2021-10-05 13:41:32.109194",
    "2021-10-05 13:41:32.109197",
    "2021-10-05 13:41:32.109200",
    "2021-10-05 13:41:32.109203",
    "2021-10-05 13:41:32.109206",
    "2021-10-05 13:41:32.109209",
    "2021-10-05 13:41:32.109212",
    "2021-10-05 13:41:32.109214",
    "2021-10-05 13:41:32.109217",
    "2021-10-05 13:41:32.109220",
    "2021-10-05 13:41:32.109222",
    "2021-10-05 13:41:32.109225",
    "2021-10-05 13:41:32.109228",
    "2021-10-05 13:41:32.109230",
    "2021-10-05 13:41:32.109234",
    "2021-10-05 13:41:32.109236",
    "2021-10-05 13:41:32.109239"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
